# #5: fund, bonds, otherbonds, antique done
# Rework the 5th sheet ("具有相當價值之財產") into the "保險" (insurance) sheet:
#  - rename the tab
#  - replace its 4-row/5-col table with a 3-row/4-col table
#  - keep the existing header-cell formatting (style used by the "index"
#    column and the first row) by copying it from a cell elsewhere in the
#    workbook that already carries that exact formatting

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("具有相當價值之財產")

$ws.Name = "保險"

# Wipe out the old 保險公司/保險名稱/要保人/備註 table entirely (values + formatting)
$ws.Range("A1:E4").Clear()

# --- New data -------------------------------------------------------------
# Row 1 (mirrors the first data row, matching the convention used by the
# other property sheets in this workbook)
$ws.Cells.Item(1, 2).Value = "富邦人壽"
$ws.Cells.Item(1, 3).Value = "生存還本保險"
$ws.Cells.Item(1, 4).Value = "林正二"

# Row 2
$ws.Cells.Item(2, 1).Value = 81
$ws.Cells.Item(2, 2).Value = "富邦人壽"
$ws.Cells.Item(2, 3).Value = "生存還本保險"
$ws.Cells.Item(2, 4).Value = "林正二"

# Row 3
$ws.Cells.Item(3, 1).Value = 82
$ws.Cells.Item(3, 2).Value = "富邦人壽"
$ws.Cells.Item(3, 3).Value = "年金保險"
$ws.Cells.Item(3, 4).Value = "林正二"

# --- Formatting -------------------------------------------------------------
# Reuse the bold/centered/bordered "header" look already present on row 1 of
# sheet 1 ("土地"), and apply it to our new row 1 and to the index column
# (column A) on the data rows - the same pattern every other sheet follows.
$srcHeader = $wb.Worksheets.Item(1).Range("B1")
$srcHeader.Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)

$srcIndex = $wb.Worksheets.Item(1).Range("A2")
$srcIndex.Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)

